# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
# This script applies a set of corrections to fixture rows in the
# "Israel Premier League" sheet. The rows below had their underlying
# match records re-mapped to the correct fixture (same matchday/date),
# so every data column except id (A), Div (C) and Date (D) is
# rewritten with the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58
$ws.Range("B58").Value = 7542500
$ws.Range("E58").Value = "Hapoel Hadera"
$ws.Range("F58").Value = "Hapoel Petah Tikva"
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = "H"
$ws.Range("L58").Value = 1.95
$ws.Range("M58").Value = 3.2
$ws.Range("N58").Value = 3.6
$ws.Range("O58").Value = 2.375
$ws.Range("P58").Value = 3
$ws.Range("Q58").Value = 2.9
$ws.Range("R58").Value = -0.25
$ws.Range("S58").Value = 2.05
$ws.Range("T58").Value = 1.8
$ws.Range("U58").Value = 2.25
$ws.Range("V58").Value = 1.925
$ws.Range("W58").Value = 1.925
$ws.Range("X58").Value = 1.375
$ws.Range("Y58").Value = -1
$ws.Range("Z58").Value = -1
$ws.Range("AA58").Value = 1.05
$ws.Range("AB58").Value = -1
$ws.Range("AC58").Value = -1
$ws.Range("AD58").Value = 0.925

# Row 59
$ws.Range("B59").Value = 7542499
$ws.Range("E59").Value = "Maccabi Petach Tikva"
$ws.Range("F59").Value = "Hapoel Beer Sheva"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 4
$ws.Range("I59").Value = 1
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = "A"
$ws.Range("L59").Value = 2.65
$ws.Range("M59").Value = 3.2
$ws.Range("N59").Value = 2.4
$ws.Range("O59").Value = 3.2
$ws.Range("P59").Value = 3.3
$ws.Range("Q59").Value = 2.05
$ws.Range("R59").Value = 0.25
$ws.Range("S59").Value = 2
$ws.Range("T59").Value = 1.85
$ws.Range("U59").Value = 2.25
$ws.Range("V59").Value = 1.85
$ws.Range("W59").Value = 2
$ws.Range("X59").Value = -1
$ws.Range("Y59").Value = -1
$ws.Range("Z59").Value = 1.05
$ws.Range("AA59").Value = -1
$ws.Range("AB59").Value = 0.8500000000000001
$ws.Range("AC59").Value = 0.8500000000000001
$ws.Range("AD59").Value = -1

# Row 60
$ws.Range("B60").Value = 7542748
$ws.Range("E60").Value = "MS Ashdod"
$ws.Range("F60").Value = "Hapoel Jerusalem FC"
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = "H"
$ws.Range("L60").Value = 2.5
$ws.Range("M60").Value = 3.2
$ws.Range("N60").Value = 2.625
$ws.Range("O60").Value = 2.4
$ws.Range("P60").Value = 2.9
$ws.Range("Q60").Value = 3
$ws.Range("R60").Value = -0.25
$ws.Range("S60").Value = 2.125
$ws.Range("T60").Value = 1.75
$ws.Range("U60").Value = 2
$ws.Range("V60").Value = 2.05
$ws.Range("W60").Value = 1.8
$ws.Range("X60").Value = 1.4
$ws.Range("Y60").Value = -1
$ws.Range("Z60").Value = -1
$ws.Range("AA60").Value = 1.125
$ws.Range("AB60").Value = -1
$ws.Range("AC60").Value = 0
$ws.Range("AD60").Value = 0

# Row 109
$ws.Range("B109").Value = 7542735
$ws.Range("E109").Value = "Hapoel Petah Tikva"
$ws.Range("F109").Value = "Maccabi Netanya"
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 2
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = "H"
$ws.Range("L109").Value = 3.75
$ws.Range("M109").Value = 3.6
$ws.Range("N109").Value = 1.909
$ws.Range("O109").Value = 3.8
$ws.Range("P109").Value = 3.75
$ws.Range("Q109").Value = 1.85
$ws.Range("R109").Value = 0.5
$ws.Range("S109").Value = 1.95
$ws.Range("T109").Value = 1.9
$ws.Range("U109").Value = 2.5
$ws.Range("V109").Value = 1.975
$ws.Range("W109").Value = 1.875
$ws.Range("X109").Value = 2.8
$ws.Range("Y109").Value = -1
$ws.Range("Z109").Value = -1
$ws.Range("AA109").Value = 0.95
$ws.Range("AB109").Value = -1
$ws.Range("AC109").Value = -1
$ws.Range("AD109").Value = 0.875

# Row 110
$ws.Range("B110").Value = 7542737
$ws.Range("E110").Value = "MS Ashdod"
$ws.Range("F110").Value = "Hapoel Haifa"
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 1
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1
$ws.Range("K110").Value = "A"
$ws.Range("L110").Value = 3
$ws.Range("M110").Value = 3.2
$ws.Range("N110").Value = 2.45
$ws.Range("O110").Value = 3.2
$ws.Range("P110").Value = 3.25
$ws.Range("Q110").Value = 2.3
$ws.Range("R110").Value = 0.25
$ws.Range("S110").Value = 1.85
$ws.Range("T110").Value = 2
$ws.Range("U110").Value = 2.25
$ws.Range("V110").Value = 1.875
$ws.Range("W110").Value = 1.975
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 1.3
$ws.Range("AA110").Value = -1
$ws.Range("AB110").Value = 1
$ws.Range("AC110").Value = -1
$ws.Range("AD110").Value = 0.9750000000000001

# Row 144
$ws.Range("B144").Value = 6799960
$ws.Range("E144").Value = "Maccabi Petach Tikva"
$ws.Range("F144").Value = "Maccabi Bnei Raina"
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 0
$ws.Range("J144").Value = 0
$ws.Range("K144").Value = "H"
$ws.Range("L144").Value = 2.625
$ws.Range("M144").Value = 3.25
$ws.Range("N144").Value = 2.5
$ws.Range("O144").Value = 2.8
$ws.Range("P144").Value = 3.25
$ws.Range("Q144").Value = 2.375
$ws.Range("R144").Value = 0.25
$ws.Range("S144").Value = 1.775
$ws.Range("T144").Value = 2.1
$ws.Range("U144").Value = 2.25
$ws.Range("V144").Value = 1.875
$ws.Range("W144").Value = 1.975
$ws.Range("X144").Value = 1.8
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.7749999999999999
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = -1
$ws.Range("AD144").Value = 0.9750000000000001

# Row 145
$ws.Range("B145").Value = 6799962
$ws.Range("E145").Value = "MS Ashdod"
$ws.Range("F145").Value = "Hapoel Petah Tikva"
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 0
$ws.Range("I145").Value = 2
$ws.Range("J145").Value = 0
$ws.Range("K145").Value = "H"
$ws.Range("L145").Value = 2.2
$ws.Range("M145").Value = 3.1
$ws.Range("N145").Value = 3.2
$ws.Range("O145").Value = 2.2
$ws.Range("P145").Value = 3.1
$ws.Range("Q145").Value = 3.2
$ws.Range("R145").Value = -0.25
$ws.Range("S145").Value = 2
$ws.Range("T145").Value = 1.85
$ws.Range("U145").Value = 2.25
$ws.Range("V145").Value = 2
$ws.Range("W145").Value = 1.85
$ws.Range("X145").Value = 1.2
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -1
$ws.Range("AA145").Value = 1
$ws.Range("AB145").Value = -1
$ws.Range("AC145").Value = -0.5
$ws.Range("AD145").Value = 0.425

# Row 220
$ws.Range("B220").Value = 8016173
$ws.Range("E220").Value = "Hapoel Hadera"
$ws.Range("F220").Value = "Maccabi Petach Tikva"
$ws.Range("G220").Value = 0
$ws.Range("H220").Value = 4
$ws.Range("I220").Value = 0
$ws.Range("J220").Value = 3
$ws.Range("K220").Value = "A"
$ws.Range("L220").Value = 2.625
$ws.Range("M220").Value = 3
$ws.Range("N220").Value = 2.7
$ws.Range("O220").Value = 2.75
$ws.Range("P220").Value = 2.7
$ws.Range("Q220").Value = 2.875
$ws.Range("R220").Value = 0
$ws.Range("S220").Value = 1.85
$ws.Range("T220").Value = 2
$ws.Range("U220").Value = 2.25
$ws.Range("V220").Value = 1.975
$ws.Range("W220").Value = 1.875
$ws.Range("X220").Value = -1
$ws.Range("Y220").Value = -1
$ws.Range("Z220").Value = 1.875
$ws.Range("AA220").Value = -1
$ws.Range("AB220").Value = 1
$ws.Range("AC220").Value = 0.9750000000000001
$ws.Range("AD220").Value = -1

# Row 221
$ws.Range("B221").Value = 8015672
$ws.Range("E221").Value = "Hapoel Bnei Sakhnin"
$ws.Range("F221").Value = "Maccabi Tel Aviv"
$ws.Range("G221").Value = 1
$ws.Range("H221").Value = 1
$ws.Range("I221").Value = 0
$ws.Range("J221").Value = 0
$ws.Range("K221").Value = "D"
$ws.Range("L221").Value = 7
$ws.Range("M221").Value = 4.333
$ws.Range("N221").Value = 1.444
$ws.Range("O221").Value = 10
$ws.Range("P221").Value = 5.25
$ws.Range("Q221").Value = 1.3
$ws.Range("R221").Value = 1.5
$ws.Range("S221").Value = 1.975
$ws.Range("T221").Value = 1.875
$ws.Range("U221").Value = 3
$ws.Range("V221").Value = 1.975
$ws.Range("W221").Value = 1.875
$ws.Range("X221").Value = -1
$ws.Range("Y221").Value = 4.25
$ws.Range("Z221").Value = -1
$ws.Range("AA221").Value = 0.9750000000000001
$ws.Range("AB221").Value = -1
$ws.Range("AC221").Value = -1
$ws.Range("AD221").Value = 0.875

# Row 236
$ws.Range("B236").Value = 8016163
$ws.Range("E236").Value = "Hapoel Haifa"
$ws.Range("F236").Value = "Maccabi Tel Aviv"
$ws.Range("G236").Value = 0
$ws.Range("H236").Value = 3
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = 0
$ws.Range("K236").Value = "A"
$ws.Range("L236").Value = 3.4
$ws.Range("M236").Value = 4
$ws.Range("N236").Value = 1.75
$ws.Range("O236").Value = 4
$ws.Range("P236").Value = 4.333
$ws.Range("Q236").Value = 1.571
$ws.Range("R236").Value = 1
$ws.Range("S236").Value = 1.8
$ws.Range("T236").Value = 2.05
$ws.Range("U236").Value = 2.75
$ws.Range("V236").Value = 1.825
$ws.Range("W236").Value = 2.025
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.571
$ws.Range("AA236").Value = -1
$ws.Range("AB236").Value = 1.05
$ws.Range("AC236").Value = 0.4125
$ws.Range("AD236").Value = -0.5

# Row 237
$ws.Range("B237").Value = 8016030
$ws.Range("E237").Value = "Hapoel Beer Sheva"
$ws.Range("F237").Value = "Maccabi Haifa"
$ws.Range("G237").Value = 1
$ws.Range("H237").Value = 4
$ws.Range("I237").Value = 1
$ws.Range("J237").Value = 3
$ws.Range("K237").Value = "A"
$ws.Range("L237").Value = 2.625
$ws.Range("M237").Value = 3.4
$ws.Range("N237").Value = 2.25
$ws.Range("O237").Value = 4
$ws.Range("P237").Value = 3.8
$ws.Range("Q237").Value = 1.727
$ws.Range("R237").Value = 0.75
$ws.Range("S237").Value = 1.85
$ws.Range("T237").Value = 2
$ws.Range("U237").Value = 2.5
$ws.Range("V237").Value = 1.825
$ws.Range("W237").Value = 2.025
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.7270000000000001
$ws.Range("AA237").Value = -1
$ws.Range("AB237").Value = 1
$ws.Range("AC237").Value = 0.825
$ws.Range("AD237").Value = -1
